$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "As an employee, I can view a customer's bank accounts" task as DONE
$ws.Range("C8").Value = "DONE"

# Update selection to reflect the newly edited cell
$ws.Range("C8").Select()
